# Insert a new data row at row 54 (pushing existing rows 54:147 down to 55:148)
# and populate it with the new weekly price-record for Orégano.
#
# This corresponds to the diff where the sheet's dimension grows from
# A1:R147 to A1:R148 and a brand-new observation (date 2022-01-20) is
# inserted ahead of the existing chronologically-ordered entries, shifting
# every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 54..147 down by inserting a new row at position 54.
$ws.Rows.Item(54).Insert()

$newRow = 54

$ws.Cells.Item($newRow, 1).Value  = 6
$ws.Cells.Item($newRow, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($newRow, 3).Value  = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value  = 44581
$ws.Cells.Item($newRow, 5).Value  = 13
$ws.Cells.Item($newRow, 6).Value  = 100112029
$ws.Cells.Item($newRow, 7).Value  = "Orégano"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 34
$ws.Cells.Item($newRow, 11).Value = 8000
$ws.Cells.Item($newRow, 12).Value = 9000
$ws.Cells.Item($newRow, 13).Value = 8441
$ws.Cells.Item($newRow, 14).Value = "$/docena de atados"
$ws.Cells.Item($newRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value = 2814
$ws.Cells.Item($newRow, 17).Value = 3
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by the
# rest of the "Fecha" column.
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
